$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Chart")

$updates = @(
    @{ Row = 2; Date = "2025-11-10"; C = 66 },
    @{ Row = 3; Date = "2025-11-11"; C = 54 },
    @{ Row = 4; Date = "2025-11-12"; C = 46 },
    @{ Row = 5; Date = "2025-11-13"; C = 43 },
    @{ Row = 6; Date = "2025-11-14"; C = 40 },
    @{ Row = 7; Date = "2025-11-15"; C = 37 },
    @{ Row = 8; Date = "2025-11-16"; C = 35 },
    @{ Row = 9; Date = "2025-11-17"; C = 30 },
    @{ Row = 10; Date = "2025-11-18"; C = 29 },
    @{ Row = 11; Date = "2025-11-19"; C = 26 },
    @{ Row = 12; Date = "2025-11-20"; C = 25 },
    @{ Row = 13; Date = "2025-11-21"; C = 25 },
    @{ Row = 14; Date = "2025-11-22"; C = 26 },
    @{ Row = 15; Date = "2025-11-23"; C = 26 },
    @{ Row = 16; Date = "2025-11-24"; C = 25 },
    @{ Row = 17; Date = "2025-11-25"; C = 25 },
    @{ Row = 18; Date = "2025-11-26"; C = 27 },
    @{ Row = 19; Date = "2025-11-27"; C = 28 },
    @{ Row = 20; Date = "2025-11-28"; C = 28 },
    @{ Row = 21; Date = "2025-11-29"; C = 27 },
    @{ Row = 22; Date = "2025-11-30"; C = 27 },
    @{ Row = 23; Date = "2025-12-01"; C = 27 },
    @{ Row = 24; Date = "2025-12-02"; C = 27 },
    @{ Row = 25; Date = "2025-12-03"; C = 27 },
    @{ Row = 26; Date = "2025-12-04"; C = 26 },
    @{ Row = 27; Date = "2025-12-05"; C = 25 },
    @{ Row = 28; Date = "2025-12-06"; C = 25 },
    @{ Row = 29; Date = "2025-12-07"; C = 25 },
    @{ Row = 30; Date = "2025-12-08"; C = 26 },
    @{ Row = 31; Date = "2025-12-09"; C = 27 },
    @{ Row = 32; Date = "2025-12-10"; C = 27 },
    @{ Row = 33; Date = "2025-12-11"; C = 29 },
    @{ Row = 34; Date = "2025-12-12"; C = 29 },
    @{ Row = 35; Date = "2025-12-13"; C = 30 },
    @{ Row = 36; Date = "2025-12-14"; C = 30 },
    @{ Row = 37; Date = "2025-12-15"; C = 31 },
    @{ Row = 38; Date = "2025-12-16"; C = 31 },
    @{ Row = 39; Date = "2025-12-17"; C = 31 },
    @{ Row = 40; Date = "2025-12-18"; C = 31 },
    @{ Row = 41; Date = "2025-12-19"; C = 31 },
    @{ Row = 42; Date = "2025-12-20"; C = 32 },
    @{ Row = 43; Date = "2025-12-21"; C = 32 },
    @{ Row = 44; Date = "2025-12-22"; C = 32 },
    @{ Row = 45; Date = "2025-12-23"; C = 32 },
    @{ Row = 46; Date = "2025-12-24"; C = 30 },
    @{ Row = 47; Date = "2025-12-25"; C = 31 },
    @{ Row = 48; Date = "2025-12-26"; C = 32 },
    @{ Row = 49; Date = "2025-12-27"; C = 30 },
    @{ Row = 50; Date = "2025-12-28"; C = 28 },
    @{ Row = 51; Date = "2025-12-29"; C = 28 },
    @{ Row = 52; Date = "2025-12-30"; C = 28 },
    @{ Row = 53; Date = "2025-12-31"; C = 28 },
    @{ Row = 54; Date = "2026-01-01"; C = 29 },
    @{ Row = 55; Date = "2026-01-02"; C = 29 },
    @{ Row = 56; Date = "2026-01-03"; C = 28 },
    @{ Row = 57; Date = "2026-01-04"; C = 27 },
    @{ Row = 58; Date = "2026-01-05"; C = 27 },
    @{ Row = 59; Date = "2026-01-06"; C = 28 },
    @{ Row = 60; Date = "2026-01-07"; C = 27 },
    @{ Row = 61; Date = "2026-01-08"; C = 27 },
    @{ Row = 62; Date = "2026-01-09"; C = 27 },
    @{ Row = 63; Date = "2026-01-10"; C = 27 },
    @{ Row = 64; Date = "2026-01-11"; C = 26 },
    @{ Row = 65; Date = "2026-01-12"; C = 26 },
    @{ Row = 66; Date = "2026-01-13"; C = 27 },
    @{ Row = 67; Date = "2026-01-14"; C = 26 },
    @{ Row = 68; Date = "2026-01-15"; C = 26 },
    @{ Row = 69; Date = "2026-01-16"; C = 25 },
    @{ Row = 70; Date = "2026-01-17"; C = 25 },
    @{ Row = 71; Date = "2026-01-18"; C = 25 },
    @{ Row = 72; Date = "2026-01-19"; C = 25 },
    @{ Row = 73; Date = "2026-01-20"; C = 26 },
    @{ Row = 74; Date = "2026-01-21"; C = 25 },
    @{ Row = 75; Date = "2026-01-22"; C = 24 },
    @{ Row = 76; Date = "2026-01-23"; C = 23 },
    @{ Row = 77; Date = "2026-01-24"; C = 24 },
    @{ Row = 78; Date = "2026-01-25"; C = 24 },
    @{ Row = 79; Date = "2026-01-26"; C = 24 },
    @{ Row = 80; Date = "2026-01-27"; C = 25 },
    @{ Row = 81; Date = "2026-01-28"; C = 26 },
    @{ Row = 82; Date = "2026-01-29"; C = 27 },
    @{ Row = 83; Date = "2026-01-30"; C = 28 },
    @{ Row = 84; Date = "2026-01-31"; C = 28 },
    @{ Row = 85; Date = "2026-02-01"; C = 28 },
    @{ Row = 86; Date = "2026-02-02"; C = 28 },
    @{ Row = 87; Date = "2026-02-03"; C = 28 },
    @{ Row = 88; Date = "2026-02-04"; C = 28 },
    @{ Row = 89; Date = "2026-02-05"; C = 28 },
    @{ Row = 90; Date = "2026-02-06"; C = 28 },
    @{ Row = 91; Date = "2026-02-07"; C = 27 }
)

foreach ($item in $updates) {
    $dateCell = $ws.Cells.Item($item.Row, 1)
    $dateCell.NumberFormat = "@"
    $dateCell.Value = $item.Date
    $dateCell.ClearFormats()

    $ws.Cells.Item($item.Row, 3).Value = $item.C
}
